# Reading-group schedule update: fix a typo in the Hill Bayesian paper link
# (donwload -> download) and add a new "presentation" file link for the
# first meeting (Christian Schulz / row 2), wiring it into column E ("Files").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "donwload" typo in the Hill et al. paper link (row 7).
$ws.Range("C7").Value = '<a href = "paper/Hill_Bayesian_Nonparametric_Modeling_for_Causal_Inference.pdf" download="Hill_Bayesian_Nonparametric_Modeling_for_Causal_Inference.pdf"> Bayesian Nonparametric Modeling for Causal Inference</a>'

# New "Files" link for the first meeting row: the presentation slides.
$ws.Range("E2").Value = '<a href = "presentation/causal_macro_reading_group_WS24.pdf" download="causal_macro_reading_group_WS24.pdf">&copy;</a>'
$ws.Range("E2").WrapText = $true

# Leave the selection on the newly added cell, matching the saved workbook.
[void]$ws.Range("E2").Select()
